# PYME-4265: add column "Subscription ID" in justifications control
#
# A new column is inserted before column I ("Antivirus Quantity"), shifting
# the remaining headers ("Antivirus Quantity" .. "Domain") one column to the
# right, and the autofilter / _FilterDatabase range + selection are fixed up
# to cover the new, wider header row (A1:M1 -> A1:N1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I, pushing "Antivirus Quantity" .. "Domain" right.
$ws.Range("I1").EntireColumn.Insert() | Out-Null

# New header cell + matching column width (same as the neighbouring
# "Antivirus Quantity" column, but not bestFit since it's a brand-new column).
$ws.Range("I1").Value = "Subscription ID"
$ws.Range("I1").ColumnWidth = 17.17

# Resize the AutoFilter so it spans the new column too.
$ws.AutoFilterMode = $false
$ws.Range("A1:N1").AutoFilter() | Out-Null

# Keep the workbook-level _FilterDatabase name in sync with the new range.
$n = $wb.Names.Item(1)
$n.RefersTo = "=Data!`$A`$1:`$N`$1"

# Match the author's resulting selection on the new column.
$ws.Range("I1").Select() | Out-Null
